# Auto-generated edit script: apply numeric updates per commit diff
# "Add budget outputs with UD penalty"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Summary")
$ws.Range("B3").Value = 0.01

$ws.Range("B6").Value = 131579.911255382

$ws.Range("B8").Value = 24253065.61257719

$ws.Range("B10").Value = 2491228.976683192

$ws = $wb.Worksheets.Item("Unmet Demand")
$arr = New-Object "object[,]" 1,15
$arr[0,0] = 414.4337959369544
$arr[0,1] = 330.5757541782243
$arr[0,2] = 176.9760193775952
$arr[0,3] = 107.2955742555736
$arr[0,4] = 109.5572237694796
$arr[0,5] = 98.64091687123928
$arr[0,6] = 77.7676953375541
$arr[0,7] = 74.3656454478664
$arr[0,8] = 83.69133109099639
$arr[0,9] = 106.2781106359148
$arr[0,10] = 128.4697750236904
$arr[0,11] = 161.0018864037399
$arr[0,12] = 189.2190633734531
$arr[0,13] = 219.2920578056454
$arr[0,14] = 251.276137581582
$ws.Range("G5:U5").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 136.8785924310737
$arr[0,1] = 107.7452501129632
$arr[0,2] = 83.51510303826707
$arr[0,3] = 82.91243519753434
$arr[0,4] = 62.76629045205057
$arr[0,5] = 37.60657774285653
$arr[0,6] = 24.33271034503603
$arr[0,7] = 10.42253800004659
$arr[0,8] = 31.97882363640291
$arr[0,9] = 45.19417501179163
$arr[0,10] = 80.63453985745144
$arr[0,11] = 116.8133877002326
$arr[0,12] = 163.0473981187501
$arr[0,13] = 198.2907557613397
$arr[0,14] = 225.9107949275447
$ws.Range("G6:U6").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 167.6012020808691
$arr[0,1] = 158.7616981666879
$arr[0,2] = 143.7288091611985
$arr[0,3] = 99.43816791380756
$arr[0,4] = 83.72822537421013
$arr[0,5] = 76.93542539304551
$arr[0,6] = 77.82642397052864
$arr[0,7] = 68.03899070462725
$arr[0,8] = 83.36329197944329
$arr[0,9] = 90.5862140395771
$arr[0,10] = 118.8664697760067
$arr[0,11] = 159.7675876048201
$arr[0,12] = 217.2238431175235
$arr[0,13] = 226.2801774240348
$arr[0,14] = 286.2977687777133
$ws.Range("G7:U7").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 414.2543098065221
$arr[0,1] = 328.737591844935
$arr[0,2] = 170.0563803341062
$arr[0,3] = 92.06191329280011
$arr[0,4] = 86.72591490550931
$arr[0,5] = 70.31665934305323
$arr[0,6] = 46.25150133729363
$arr[0,7] = 42.33948647951536
$arr[0,8] = 53.44993733213093
$arr[0,9] = 80.46778072209531
$arr[0,10] = 109.0872921559756
$arr[0,11] = 149.7272407629755
$arr[0,12] = 185.1290231762283
$arr[0,13] = 218.5063572696783
$arr[0,14] = 251.2617786911475
$ws.Range("G8:U8").Value = $arr

$arr = New-Object "object[,]" 1,7
$arr[0,0] = 136.7825588849655
$arr[0,1] = 106.817768233445
$arr[0,2] = 80.20868489375378
$arr[0,3] = 73.83937108913167
$arr[0,4] = 47.25897875440141
$arr[0,5] = 16.75508344424118
$arr[0,6] = 0
$ws.Range("G9:M9").Value = $arr

$arr = New-Object "object[,]" 1,7
$arr[0,0] = 9.130000058643361
$arr[0,1] = 26.85597970277001
$arr[0,2] = 68.37594193669796
$arr[0,3] = 110.8508838462542
$arr[0,4] = 161.2636171197674
$arr[0,5] = 197.9036731785617
$arr[0,6] = 225.9044769310903
$ws.Range("O9:U9").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 167.5206907817879
$arr[0,1] = 158.0458795257661
$arr[0,2] = 141.3076148215571
$arr[0,3] = 93.74601906876785
$arr[0,4] = 74.37427626277801
$arr[0,5] = 64.96559098237566
$arr[0,6] = 65.20591187910269
$arr[0,7] = 55.71856618250388
$arr[0,8] = 71.98338581476828
$arr[0,9] = 80.84873837615825
$arr[0,10] = 112.1247463593087
$arr[0,11] = 156.1475068297698
$arr[0,12] = 215.8207507508087
$arr[0,13] = 225.9361746006879
$arr[0,14] = 286.2933772523089
$ws.Range("G10:U10").Value = $arr

$ws = $wb.Worksheets.Item("Household Surplus")
$ws.Range("B3").Value = 168381.7769200939

$ws.Range("B4").Value = 202375.4442256468

$ws = $wb.Worksheets.Item("Costs and Revenues")
$arr = New-Object "object[,]" 1,2
$arr[0,0] = 62730.46591140758
$arr[0,1] = 75394.77333896644
$ws.Range("C2:D2").Value = $arr

$ws.Range("H2").Value = 91976.24205358137

$ws.Range("K2").Value = 91976.24205358134

$ws.Range("P2").Value = 91976.24205358137

$arr = New-Object "object[,]" 1,3
$arr[0,0] = 196825.9098199031
$arr[0,1] = 38236.46568336456
$arr[0,2] = 52530.53686621619
$ws.Range("C3:E3").Value = $arr

$arr = New-Object "object[,]" 1,2
$arr[0,0] = 38339.65294307929
$arr[0,1] = 39312.96135688073
$ws.Range("C5:D5").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = -54153.64424660708
$arr[0,1] = -189659.5376291605
$arr[0,2] = -18712.55198267718
$arr[0,3] = 16649.76056833225
$arr[0,4] = 69180.29743454847
$arr[0,5] = 69180.29743454847
$arr[0,6] = 69180.29743454847
$arr[0,7] = 69180.29743454847
$arr[0,8] = 69180.29743454847
$arr[0,9] = 69180.29743454844
$arr[0,10] = 69180.29743454844
$arr[0,11] = 69180.29743454844
$arr[0,12] = 69180.29743454844
$arr[0,13] = 69180.29743454844
$arr[0,14] = 69180.29743454847
$ws.Range("B6:P6").Value = $arr

$ws = $wb.Worksheets.Item("Installed Capacities")
$arr = New-Object "object[,]" 1,2
$arr[0,0] = 216.1492175724446
$arr[0,1] = 260.7963925174648
$ws.Range("C3:D3").Value = $arr

$ws = $wb.Worksheets.Item("Added Capacities")
$arr = New-Object "object[,]" 1,3
$arr[0,0] = 216.1492175724445
$arr[0,1] = 44.64717494502023
$arr[0,2] = 65.38503947111997
$ws.Range("C3:E3").Value = $arr

$ws = $wb.Worksheets.Item("PV Dispatch")
$arr = New-Object "object[,]" 1,15
$arr[0,0] = 0.8689415781806812
$arr[0,1] = 8.899047937542903
$arr[0,2] = 33.49987019281074
$arr[0,3] = 73.75033027111266
$arr[0,4] = 110.5326272755009
$arr[0,5] = 137.125498098748
$arr[0,6] = 152.5785378897186
$arr[0,7] = 155.0474181487245
$arr[0,8] = 146.4068803306903
$arr[0,9] = 124.9548851193548
$arr[0,10] = 93.8359148507591
$arr[0,11] = 54.58365141039226
$arr[0,12] = 19.80100621279229
$arr[0,13] = 3.803791758485934
$arr[0,14] = 0.06951532625445447
$ws.Range("G5:U5").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 0.4649247321369563
$arr[0,1] = 4.490194123533237
$arr[0,2] = 16.00727696173293
$arr[0,3] = 43.92519146913236
$arr[0,4] = 75.07514852230842
$arr[0,5] = 100.9478020370177
$arr[0,6] = 117.8013235769823
$arr[0,7] = 120.9191740832867
$arr[0,8] = 110.6174208080415
$arr[0,9] = 88.78023240253862
$arr[0,10] = 59.34723422857008
$arr[0,11] = 28.86611626373139
$arr[0,12] = 8.635772985087758
$arr[0,13] = 1.873972933481854
$arr[0,14] = 0.03058715343006293
$ws.Range("G6:U6").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 0.3897772775896541
$arr[0,1] = 3.465474340751655
$arr[0,2] = 11.72166576605978
$arr[0,3] = 27.55725352558855
$arr[0,4] = 45.28503279632526
$arr[0,5] = 57.94925088819277
$arr[0,6] = 61.09935997707642
$arr[0,7] = 59.64655376060593
$arr[0,8] = 55.0932464723995
$arr[0,9] = 47.14179000957051
$arr[0,10] = 32.63853185343913
$arr[0,11] = 17.52580377234936
$arr[0,12] = 6.792754919448789
$arr[0,13] = 1.665412004246704
$arr[0,14] = 0.02126057877761752
$ws.Range("G7:U7").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 1.048427708612923
$arr[0,1] = 10.7372102708321
$arr[0,2] = 40.41950923629976
$arr[0,3] = 88.98399123388617
$arr[0,4] = 133.3639361394712
$arr[0,5] = 165.449755626934
$arr[0,6] = 184.0947318899791
$arr[0,7] = 187.0735771170756
$arr[0,8] = 176.6482740895558
$arr[0,9] = 150.7652150331742
$arr[0,10] = 113.2183977184739
$arr[0,11] = 65.8582970511566
$arr[0,12] = 23.89104641001701
$arr[0,13] = 4.589492294453073
$arr[0,14] = 0.08387421668903385
$ws.Range("G8:U8").Value = $arr

$arr = New-Object "object[,]" 1,7
$arr[0,0] = 0.560958278245113
$arr[0,1] = 5.417676003051488
$arr[0,2] = 19.31369510624622
$arr[0,3] = 52.99825557753501
$arr[0,4] = 90.58246021995758
$arr[0,5] = 121.799296335633
$arr[0,6] = 142.1340339220183
$ws.Range("G9:M9").Value = $arr

$arr = New-Object "object[,]" 1,7
$arr[0,0] = 133.4662443858011
$arr[0,1] = 107.1184277115602
$arr[0,2] = 71.60583214932356
$arr[0,3] = 34.82862011770975
$arr[0,4] = 10.41955398407041
$arr[0,5] = 2.261055516259907
$arr[0,6] = 0.03690514988454693
$ws.Range("O9:U9").Value = $arr

$arr = New-Object "object[,]" 1,15
$arr[0,0] = 0.4702885766708382
$arr[0,1] = 4.181292981673455
$arr[0,2] = 14.14286010570121
$arr[0,3] = 33.24940237062826
$arr[0,4] = 54.63898190775737
$arr[0,5] = 69.91908529886263
$arr[0,6] = 73.71987206850237
$arr[0,7] = 71.96697828272931
$arr[0,8] = 66.47315263707451
$arr[0,9] = 56.87926567298936
$arr[0,10] = 39.38025527013719
$arr[0,11] = 21.14588454739968
$arr[0,12] = 8.195847286163604
$arr[0,13] = 2.009414827593581
$arr[0,14] = 0.02565210418204575
$ws.Range("G10:U10").Value = $arr

Write-Host "Applied all updates"